$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 83, pushing the existing rows 83-139 down to 84-140
$ws.Rows("83:83").Insert()

# Populate the new row 83 with the new data record
$ws.Range("A83").Value = 7
$ws.Range("B83").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C83").Value = "Ñuble"
$ws.Range("D83").Value = 44574
$ws.Range("E83").Value = 16
$ws.Range("F83").Value = 100112024
$ws.Range("G83").Value = "Choclo"
$ws.Range("H83").Value = "Choclero"
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 30000
$ws.Range("K83").Value = 200
$ws.Range("L83").Value = 250
$ws.Range("M83").Value = 225
$ws.Range("N83").Value = "$/unidad"
$ws.Range("O83").Value = "Región del Maule"
$ws.Range("P83").Value = 225
$ws.Range("Q83").Value = 1
$ws.Range("R83").Value = "Hortaliza"
